# Updated main GSC export data.
# The GSC export window rolled forward by one day: the oldest day
# (2025-10-05) dropped off the front of the report, so that row is
# removed from the "Chart" sheet and every following row shifts up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the entire second row (the oldest date, 2025-10-05); this
# shifts all subsequent rows up by one and shrinks the used range from
# A1:D90 down to A1:D89, matching the refreshed export.
$ws.Rows.Item(2).Delete()
